# Update (Analyze PO & Forecast)
# The forecast run was regenerated ~4 weeks later: Week_Start_Date values on the
# "Forecast Comparison" sheet shift back by 4 weeks, MyForecast (col D) is
# recomputed for the first 9 weeks, and the stale is_holiday_week flag (col J,
# previously always FALSE) is cleared. The "Summary" sheet's forecast-derived
# metrics are refreshed to match.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New Week_Start_Date (col B) and MyForecast (col D, only for rows whose
# forecast value changed) for each data row of the Forecast Comparison sheet.
$weekUpdates = @(
    @{ Row = 2;  Date = "2024-12-29"; Forecast = 274 },
    @{ Row = 3;  Date = "2025-01-05"; Forecast = 302 },
    @{ Row = 4;  Date = "2025-01-12"; Forecast = 272 },
    @{ Row = 5;  Date = "2025-01-19"; Forecast = 293 },
    @{ Row = 6;  Date = "2025-01-26"; Forecast = 290 },
    @{ Row = 7;  Date = "2025-02-02"; Forecast = 285 },
    @{ Row = 8;  Date = "2025-02-09"; Forecast = 282 },
    @{ Row = 9;  Date = "2025-02-16"; Forecast = 283 },
    @{ Row = 10; Date = "2025-02-23"; Forecast = 283 },
    @{ Row = 11; Date = "2025-03-02"; Forecast = $null },
    @{ Row = 12; Date = "2025-03-09"; Forecast = $null },
    @{ Row = 13; Date = "2025-03-16"; Forecast = $null },
    @{ Row = 14; Date = "2025-03-23"; Forecast = $null },
    @{ Row = 15; Date = "2025-03-30"; Forecast = $null },
    @{ Row = 16; Date = "2025-04-06"; Forecast = $null },
    @{ Row = 17; Date = "2025-04-13"; Forecast = $null }
)

foreach ($u in $weekUpdates) {
    $r = $u.Row

    # Column B: Week_Start_Date - keep it stored as literal text (matches the
    # original inline-string dates rather than letting Excel coerce it to a
    # date serial number).
    $cellB = $wsForecast.Cells.Item($r, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $u.Date

    # Column D: MyForecast - only rewritten where the recomputed value differs.
    if ($null -ne $u.Forecast) {
        $wsForecast.Cells.Item($r, 4).Value = $u.Forecast
    }

    # Column J: is_holiday_week - was a stale FALSE for every row; now blank.
    $wsForecast.Cells.Item($r, 10).ClearContents()
}

# Summary sheet: forecast-derived metrics recalculated against the refreshed
# Forecast Comparison data above. All "Value" cells on this sheet are stored
# as literal text (even the numeric-looking ones), so force text formatting
# before assigning each one - otherwise Excel would coerce them to numbers
# or dates and change the underlying cell type.
$summaryUpdates = @(
    @{ Row = 9;  Value = "4486" },
    @{ Row = 10; Value = "2280" },
    @{ Row = 11; Value = "1141" },
    @{ Row = 12; Value = "302" },
    @{ Row = 13; Value = "2025-01-05" },
    @{ Row = 14; Value = "271" },
    @{ Row = 15; Value = "2025-03-30" }
)

foreach ($u in $summaryUpdates) {
    $cell = $wsSummary.Cells.Item($u.Row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
